$p = $ppt.ActivePresentation
